# Payment_records.xlsx update — "Update README.md + Resources"
#
# Adds a 5% rate (column M) to several existing rows whose Screws
# total (column N = L*M) had been left at 0, fills in four new
# purchase rows (40-43, shifted to 41-44 in the sheet) with date,
# material, quantity/price columns, and nudges the saved view
# (scrolled position + active cell selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column M (rate) back-fill for rows whose N formula (=L*M) was
#     still resolving to 0 because M was blank -------------------------
$rateRows = 13, 14, 19, 20, 21, 22, 23, 24, 25, 33
foreach ($r in $rateRows) {
    $ws.Cells.Item($r, 13).Value = 0.05   # column M
}

# --- New rows of purchases -------------------------------------------
# Row 40: date + material (new "Screws and Nuts" string) + P/Q pricing
$ws.Cells.Item(40, 2).Value = 44938          # B40 (date serial)
$ws.Cells.Item(40, 3).Value = "Screws and Nuts"  # C40
$ws.Cells.Item(40, 16).Value = 38            # P40
$ws.Cells.Item(40, 17).Value = 0.05          # Q40

# Row 41: date + material + L/M pricing
$ws.Cells.Item(41, 2).Value = 44938          # B41
$ws.Cells.Item(41, 3).Value = "PETG"         # C41
$ws.Cells.Item(41, 12).Value = 48            # L41
$ws.Cells.Item(41, 13).Value = 0.05          # M41

# Row 42: date + material + P/Q pricing
$ws.Cells.Item(42, 2).Value = 44939          # B42
$ws.Cells.Item(42, 3).Value = "Screws"       # C42
$ws.Cells.Item(42, 16).Value = 12            # P42
$ws.Cells.Item(42, 17).Value = 0.05          # Q42

# Row 43: date + material + P/Q pricing
$ws.Cells.Item(43, 2).Value = 44939          # B43
$ws.Cells.Item(43, 3).Value = "Screws"       # C43
$ws.Cells.Item(43, 16).Value = 10            # P43
$ws.Cells.Item(43, 17).Value = 0.05          # Q43

# Row 44: date + material + P/Q pricing
$ws.Cells.Item(44, 2).Value = 44939          # B44
$ws.Cells.Item(44, 3).Value = "Screws"       # C44
$ws.Cells.Item(44, 16).Value = 7.4           # P44
$ws.Cells.Item(44, 17).Value = 0.05          # Q44

# --- Apply the same date number format used elsewhere in column B ----
$ws.Range("B40:B44").NumberFormat = $ws.Range("B39").NumberFormat

# --- Scroll/selection state, matching the saved view ------------------
$ws.Range("M23").Select()
$excel.ActiveWindow.ScrollRow = 9
